$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9047.200000000001
$ws.Range("I74").Value = 9047.200000000001
$ws.Range("K74").Value = 9047.200000000001
$ws.Range("M74").Value = -8111.200000000001
$ws.Range("H77").Value = 9047.200000000001
$ws.Range("I77").Value = 9047.200000000001
$ws.Range("K77").Value = 45236
$ws.Range("M77").Value = -40556
$ws.Range("H98").Value = 1717
$ws.Range("I98").Value = 1353.3846
$ws.Range("J98").Value = 4868.3335
$ws.Range("K98").Value = 1353.3846
$ws.Range("L98").Value = 4868.3335
$ws.Range("M98").Value = 144.6153999999999
$ws.Range("N98").Value = -7864.3335
$ws.Range("H100").Value = 5097.24
$ws.Range("J100").Value = 7555.2856
$ws.Range("L100").Value = 7555.2856
$ws.Range("N100").Value = -8637.285599999999
$ws.Range("H122").Value = 1717
$ws.Range("I122").Value = 1353.3846
$ws.Range("J122").Value = 4868.3335
$ws.Range("K122").Value = 4060.1538
$ws.Range("L122").Value = 14605.0005
$ws.Range("M122").Value = -1610.1538
$ws.Range("N122").Value = -19505.0005
$ws.Range("H129").Value = 3051.25
$ws.Range("J129").Value = 4370.3335
$ws.Range("L129").Value = 13111.0005
$ws.Range("N129").Value = -23111.0005
$ws.Range("H132").Value = 4688.225
$ws.Range("J132").Value = 13185
$ws.Range("L132").Value = 39555
$ws.Range("N132").Value = -44615
$ws.Range("H135").Value = 3736.6216
$ws.Range("I135").Value = 3366.6155
$ws.Range("J135").Value = 4611.1816
$ws.Range("K135").Value = 30299.5395
$ws.Range("L135").Value = 41500.6344
$ws.Range("M135").Value = -27764.5395
$ws.Range("N135").Value = -46570.6344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3364.25
$ws.Range("I2").Value = 2862.1
$ws.Range("K2").Value = 2862.1
$ws.Range("M2").Value = -2749.1
$ws.Range("H110").Value = 5765.923
$ws.Range("I110").Value = 6817.5557
$ws.Range("K110").Value = 6817.5557
$ws.Range("M110").Value = -4772.5557
$ws.Range("H116").Value = 3364.25
$ws.Range("I116").Value = 2862.1
$ws.Range("K116").Value = 2862.1
$ws.Range("M116").Value = -568.0999999999999
$ws.Range("H123").Value = 330000
$ws.Range("J123").Value = 330000
$ws.Range("L123").Value = 330000
$ws.Range("N123").Value = -339800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3364.25
$ws.Range("I3").Value = 2862.1
$ws.Range("K3").Value = 2862.1
$ws.Range("M3").Value = -2748.1
$ws.Range("H94").Value = 2058.4688
$ws.Range("I94").Value = 2263.1667
$ws.Range("K94").Value = 2263.1667
$ws.Range("M94").Value = -1812.1667
$ws.Range("H99").Value = 3083.0527
$ws.Range("I99").Value = 3056.375
$ws.Range("J99").Value = 3225.3333
$ws.Range("K99").Value = 3056.375
$ws.Range("L99").Value = 3225.3333
$ws.Range("M99").Value = -1558.375
$ws.Range("N99").Value = -6221.3333
$ws.Range("H134").Value = 2668.6086
$ws.Range("I134").Value = 2139.1177
$ws.Range("J134").Value = 4168.8335
$ws.Range("K134").Value = 6417.353099999999
$ws.Range("L134").Value = 12506.5005
$ws.Range("M134").Value = -3882.353099999999
$ws.Range("N134").Value = -17576.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15271.353
$ws.Range("I99").Value = 4799.9287
$ws.Range("K99").Value = 4799.9287
$ws.Range("M99").Value = -3301.9287
$ws.Range("H126").Value = 15271.353
$ws.Range("I126").Value = 4799.9287
$ws.Range("K126").Value = 14399.7861
$ws.Range("M126").Value = -11929.7861
$ws.Range("H132").Value = 2299.6667
$ws.Range("I132").Value = 2299.6667
$ws.Range("K132").Value = 6899.000100000001
$ws.Range("M132").Value = -4369.000100000001
$ws.Range("H134").Value = 5221.857
$ws.Range("I134").Value = 4996.2173
$ws.Range("J134").Value = 6259.8
$ws.Range("K134").Value = 14988.6519
$ws.Range("L134").Value = 18779.4
$ws.Range("M134").Value = -12453.6519
$ws.Range("N134").Value = -23849.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1562.2222
$ws.Range("I5").Value = 1381.8
$ws.Range("J5").Value = 1787.75
$ws.Range("K5").Value = 4145.4
$ws.Range("L5").Value = 5363.25
$ws.Range("M5").Value = -4033.4
$ws.Range("N5").Value = -5587.25
$ws.Range("H68").Value = 2362.6
$ws.Range("I68").Value = 1324.75
$ws.Range("J68").Value = 2740
$ws.Range("K68").Value = 3974.25
$ws.Range("L68").Value = 8220
$ws.Range("M68").Value = -3163.25
$ws.Range("N68").Value = -9842
$ws.Range("H71").Value = 2362.6
$ws.Range("I71").Value = 1324.75
$ws.Range("J71").Value = 2740
$ws.Range("K71").Value = 11922.75
$ws.Range("L71").Value = 24660
$ws.Range("M71").Value = -7866.75
$ws.Range("N71").Value = -32772
$ws.Range("H86").Value = 2947
$ws.Range("J86").Value = 3429.3333
$ws.Range("L86").Value = 10287.9999
$ws.Range("N86").Value = -12659.9999
$ws.Range("H89").Value = 2947
$ws.Range("J89").Value = 3429.3333
$ws.Range("L89").Value = 30863.9997
$ws.Range("N89").Value = -42719.9997
$ws.Range("H132").Value = 2199.4
$ws.Range("J132").Value = 2999.6667
$ws.Range("L132").Value = 26997.0003
$ws.Range("N132").Value = -32057.0003
$ws.Range("H135").Value = 1562.2222
$ws.Range("I135").Value = 1381.8
$ws.Range("J135").Value = 1787.75
$ws.Range("K135").Value = 12436.2
$ws.Range("L135").Value = 16089.75
$ws.Range("M135").Value = -9901.199999999999
$ws.Range("N135").Value = -21159.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2725.6155
$ws.Range("J122").Value = 3164.1667
$ws.Range("L122").Value = 9492.500100000001
$ws.Range("N122").Value = -14392.5001
$ws.Range("H132").Value = 2154.3728
$ws.Range("I132").Value = 2139.3542
$ws.Range("J132").Value = 2219.9092
$ws.Range("K132").Value = 6418.062600000001
$ws.Range("L132").Value = 6659.7276
$ws.Range("M132").Value = -3888.062600000001
$ws.Range("N132").Value = -11719.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 955.1818
$ws.Range("I55").Value = 452.54544
$ws.Range("J55").Value = 1457.8182
$ws.Range("K55").Value = 452.54544
$ws.Range("L55").Value = 1457.8182
$ws.Range("M55").Value = -279.54544
$ws.Range("N55").Value = -1803.8182
$ws.Range("H103").Value = 50230.2
$ws.Range("J103").Value = 50230.2
$ws.Range("L103").Value = 50230.2
$ws.Range("N103").Value = -52574.2
$ws.Range("H132").Value = 3468.8635
$ws.Range("I132").Value = 2230.75
$ws.Range("J132").Value = 4954.6
$ws.Range("K132").Value = 6692.25
$ws.Range("L132").Value = 14863.8
$ws.Range("M132").Value = -4162.25
$ws.Range("N132").Value = -19923.8
$ws.Range("H136").Value = 7742.278
$ws.Range("I136").Value = 8364.532999999999
$ws.Range("J136").Value = 4631
$ws.Range("K136").Value = 25093.599
$ws.Range("L136").Value = 13893
$ws.Range("M136").Value = -22543.599
$ws.Range("N136").Value = -18993

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8077.8335
$ws.Range("J41").Value = 8077.8335
$ws.Range("L41").Value = 8077.8335
$ws.Range("N41").Value = -8857.833500000001
$ws.Range("H126").Value = 4252.2085
$ws.Range("I126").Value = 4432.4116
$ws.Range("J126").Value = 3814.5715
$ws.Range("K126").Value = 13297.2348
$ws.Range("L126").Value = 11443.7145
$ws.Range("M126").Value = -10827.2348
$ws.Range("N126").Value = -16383.7145
$ws.Range("H132").Value = 1531.3721
$ws.Range("I132").Value = 1357.0526
$ws.Range("J132").Value = 2856.2
$ws.Range("K132").Value = 4071.1578
$ws.Range("L132").Value = 8568.599999999999
$ws.Range("M132").Value = -1541.1578
$ws.Range("N132").Value = -13628.6
$ws.Range("H136").Value = 11838
$ws.Range("I136").Value = 12022.06
$ws.Range("J136").Value = 10687.625
$ws.Range("K136").Value = 36066.18
$ws.Range("L136").Value = 32062.875
$ws.Range("M136").Value = -33516.18
$ws.Range("N136").Value = -37162.875
